$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 2 to hold the "slug" identifiers that relate
# each data column to its header, shifting the existing rows 2-4 down to 3-5.
$ws.Rows.Item(2).Insert()

# The inserted row picks up a generic default format; copy the header row's
# formatting onto it so it matches the rest of the sheet (style "1").
$ws.Range("A1:F1").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)

$ws.Range("A2").Value = "division-2-digitos-descripcion"
$ws.Range("B2").Value = "afiliaciones-en-alta"
$ws.Range("C2").Value = "comarca-nombre"
$ws.Range("D2").Value = "comarca-codigo"
$ws.Range("E2").Value = "division-2-digitos-codigo"
$ws.Range("F2").Value = "mes-y-ano"
